$wb = $excel.ActiveWorkbook

# --- Sheet "RVL" (sheet1) ---
$ws1 = $wb.Worksheets.Item("RVL")

# Remove the 4 "Param" rows (fromRow/fromCol/toRow/toCol) that followed the
# Database "Range" action block (old rows 9-12).
$ws1.Rows.Item(9).Resize(4).Delete()

# The "Functions" object block (now rows 11-14) is renamed to "Sfdc", and its
# actions drop the redundant "Sfdc" prefix.
$ws1.Range("C11").Value = "Sfdc"
$ws1.Range("D11").Value = "Launch"

$ws1.Range("C12").Value = "Sfdc"
$ws1.Range("D12").Value = "OpenApp"

$ws1.Range("C13").Value = "Sfdc"
$ws1.Range("D13").Value = "NavigateModule"

$ws1.Range("C14").Value = "Sfdc"
$ws1.Range("D14").Value = "SelectListView"

# --- Sheet "Cleanup" (sheet2) ---
$ws2 = $wb.Worksheets.Item("Cleanup")

# Remove two blank formatting-only rows (old rows 3-4).
$ws2.Rows.Item(3).Resize(2).Delete()
